$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D:E value cells to be stored as text (they hold formatted
# numeric-looking strings like "5.28" or "54.017.69" that Excel
# would otherwise auto-convert to numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '54.017.69'
$ws.Range("E2").Value = '  -6.92%  '
$ws.Range("D3").Value = '2.409.80'
$ws.Range("E3").Value = '  -10.07%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '465.61'
$ws.Range("E5").Value = '  -6.35%  '
$ws.Range("D6").Value = '130.03'
$ws.Range("E6").Value = '  -5.53%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '0.490'
$ws.Range("E8").Value = '  -6.49%  '
$ws.Range("D9").Value = '2.418.28'
$ws.Range("E9").Value = '  -9.89%  '
$ws.Range("D10").Value = '0.0946'
$ws.Range("E10").Value = '  -8.90%  '
$ws.Range("D11").Value = '5.28'
$ws.Range("E11").Value = '  -12.07%  '
$ws.Range("E12").Value = '  -8.75%  '
$ws.Range("E13").Value = '  -3.83%  '
$ws.Range("D14").Value = '2.842.90'
$ws.Range("E14").Value = '  -9.88%  '
$ws.Range("D15").Value = '54.203.27'
$ws.Range("E15").Value = '  -6.78%  '
$ws.Range("D16").Value = '0.0000133'
$ws.Range("D17").Value = '19.56'
$ws.Range("E17").Value = '  -7.88%  '
$ws.Range("D18").Value = '2.437.76'
$ws.Range("E18").Value = '  -9.72%  '
$ws.Range("D19").Value = '4.21'
$ws.Range("E19").Value = '  -10.02%  '
$ws.Range("D20").Value = '310.19'
$ws.Range("E20").Value = '  -6.47%  '
$ws.Range("D21").Value = '9.51'
$ws.Range("E21").Value = '  -12.58%  '
$ws.Range("D22").Value = '0.994'
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").Value = '5.67'
$ws.Range("E23").Value = '  +0.50%  '
$ws.Range("D24").Value = '5.38'
$ws.Range("E24").Value = '  -12.80%  '
$ws.Range("D25").Value = '56.18'
$ws.Range("E25").Value = '  -9.96%  '
$ws.Range("E26").Value = '  +1.11%  '
$ws.Range("D27").Value = '0.385'
$ws.Range("E27").Value = '  -8.64%  '
$ws.Range("D28").Value = '2.547.79'
$ws.Range("E28").Value = '  -9.41%  '
$ws.Range("E29").Value = '  -7.65%  '
$ws.Range("D30").Value = '7.08'
$ws.Range("E30").Value = '  -3.38%  '
$ws.Range("E31").Value = '  +0.09%  '
$ws.Range("D32").Value = '0.0₃0708'
$ws.Range("E32").Value = '  -12.36%  '
$ws.Range("D33").Value = '145.73'
$ws.Range("E33").Value = '  -2.78%  '
$ws.Range("D34").Value = '17.71'
$ws.Range("E34").Value = '  -6.52%  '
$ws.Range("E35").Value = '  -9.56%  '
$ws.Range("D36").Value = '5.00'
$ws.Range("E36").Value = '  -6.12%  '
$ws.Range("D37").Value = '3.53'
$ws.Range("E37").Value = '  -14.76%  '
$ws.Range("E38").Value = '  -5.06%  '
$ws.Range("E39").Value = '  -14.28%  '
$ws.Range("D40").Value = '0.997'
$ws.Range("E40").Value = '  -0.09%  '
$ws.Range("D41").Value = '32.64'
$ws.Range("E41").Value = '  -7.13%  '
$ws.Range("D42").Value = '0.593'
$ws.Range("E42").Value = '  +1.22%  '
$ws.Range("D43").Value = '0.0522'
$ws.Range("E43").Value = '  -5.38%  '
$ws.Range("E44").Value = '  -7.92%  '
$ws.Range("D45").Value = '10.10'
$ws.Range("E45").Value = '  -2.44%  '
$ws.Range("E46").Value = '  -9.63%  '
$ws.Range("D47").Value = '1.919.45'
$ws.Range("E47").Value = '  -11.42%  '
$ws.Range("D48").Value = '0.0881'
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("E49").Value = '  -3.25%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '16.51'
$ws.Range("E50").Value = '  -10.74%  '
$ws.Range("B51").Value = 'Bittensor'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D51").Value = '230.01'
$ws.Range("E51").Value = '  +5.82%  '

# Restore the default (unstyled) cell style so no stray number format
# is left attached to the cells (matches original workbook styling).
$ws.Range("D2:E51").Style = "Normal"
